# postproc and writing added
#
# 1. Add a new "setup" worksheet at the end of the workbook containing
#    calorimeter setup info (Calorimeter/DSC, Initial volume/15).
# 2. Update the log10(K) input value on "input_k_constants_log10" from
#    5.13 to 4, and make that sheet the active/selected one.
# 3. "heats" is no longer the selected sheet (its selection otherwise
#    stays the same).

$wb = $excel.ActiveWorkbook

# --- 1. Add the new "setup" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "setup"

$newSheet.Range("A1").Value = "Calorimeter"
$newSheet.Range("B1").Value = "DSC"
$newSheet.Range("A2").Value = "Initial volume"
$newSheet.Range("B2").Value = 15

$newSheet.Range("A2").Select() | Out-Null

# --- 2. Update the log K value and make this sheet the active tab ---
$ws2 = $wb.Worksheets.Item("input_k_constants_log10")
$ws2.Range("A2").Value = 4

$ws2.Activate() | Out-Null
$ws2.Range("A3").Select() | Out-Null
